$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Start clean: wipe all existing cell content/formatting on the sheet ---
$ws.Cells.Clear()

# --- Row 4: headers ---
$ws.Range("A4").Value2 = "lane"
$ws.Range("C4").Value2 = "auto"
$ws.Range("D4").Value2 = "norm"
$ws.Range("F4").Value2 = "timeRed"
$ws.Range("G4").Value2 = "norm"
$ws.Range("H4").Value2 = "bonus"
$ws.Range("J4").Value2 = "time green"
$ws.Range("K4").Value2 = "norm"
$ws.Range("O4").Value2 = "out"

# --- Row 1 / Row 2: outIndex lookup helper ---
$ws.Range("O1").Value2 = "outIndex"
$ws.Range("O2").Formula = "=MATCH(O4,A4:R4)"

# --- Column A: lane numbers ---
$ws.Range("A5").Value2 = 1
$ws.Range("A6").Value2 = 2
$ws.Range("A7").Value2 = 3
$ws.Range("A8").Value2 = 4
$ws.Range("A9").Value2 = 5
$ws.Range("A10").Value2 = 6

# --- Column C: auto (raw) ---
$ws.Range("C5").Value2 = 1
$ws.Range("C6").Value2 = 1
$ws.Range("C7").Value2 = 2
$ws.Range("C8").Value2 = 2
$ws.Range("C9").Value2 = 1
$ws.Range("C10").Value2 = 1

# --- Column D: norm(auto) - D5 individual, D6:D10 shared ---
$ws.Range("D5").Formula = "=(C5 - `$D`$12) / (`$D`$13 - `$D`$12)"
$ws.Range("D6:D10").Formula = "=(C6 - `$D`$12) / (`$D`$13 - `$D`$12)"

# --- Column F: timeRed (raw) ---
$ws.Range("F5").Value2 = 40
$ws.Range("F6").Value2 = 50
$ws.Range("F7").Value2 = 0
$ws.Range("F8").Value2 = 0
$ws.Range("F9").Value2 = 56
$ws.Range("F10").Value2 = 362

# --- Column G: norm(timeRed) - G5 individual, G6:G10 shared ---
$ws.Range("G5").Formula = "=(F5 - `$G`$12) / (`$G`$13 - `$G`$12)"
$ws.Range("G6:G10").Formula = "=(F6 - `$G`$12) / (`$G`$13 - `$G`$12)"

# --- Column H: bonus - H5 individual, H6:H9 shared, H10 individual (full abs refs) ---
$ws.Range("H5").Formula = "=IF(F5>`$G`$13,G5+(F5-`$G`$13),G5)"
$ws.Range("H6:H9").Formula = "=IF(F6>`$G`$13,G6+(F6-`$G`$13),G6)"

# H10 written individually (not part of the H6:H9 shared-formula group)
$ws.Range("H10").Formula = "=IF(F10>`$G`$13,G10+(F10-`$G`$13),G10)"

# --- Column J: time green (raw) ---
$ws.Range("J5").Value2 = 0
$ws.Range("J6").Value2 = 0
$ws.Range("J7").Value2 = -1
$ws.Range("J8").Value2 = 4
$ws.Range("J9").Value2 = 0
$ws.Range("J10").Value2 = 0

# --- Column K: norm(time green) - K5 individual, K6:K10 shared ---
$ws.Range("K5").Formula = "=(J5 - `$K`$12) / (`$K`$13 - `$K`$12)"
$ws.Range("K6:K10").Formula = "=(J6 - `$K`$12) / (`$K`$13 - `$K`$12)"

# --- Column O: out - O5 individual, O6:O10 shared ---
$ws.Range("O5").Formula = "=H5*D5+D5*K5"
$ws.Range("O6:O10").Formula = "=H6*D6+D6*K6"

# --- Row 12: min reference values ---
$ws.Range("C12").Value2 = "min"
$ws.Range("D12").Value2 = 0
$ws.Range("F12").Value2 = "min"
$ws.Range("G12").Value2 = 0
$ws.Range("J12").Value2 = "min"
$ws.Range("K12").Value2 = -0.01

# --- Row 13: max reference values ---
$ws.Range("C13").Value2 = "max"
$ws.Range("D13").Formula = "=MAX(C5:C10)"
$ws.Range("F13").Value2 = "max"
$ws.Range("G13").Value2 = 360
$ws.Range("J13").Value2 = "max"
$ws.Range("K13").Value2 = 180

# --- Cosmetic: column width / view state ---
$ws.Columns("P").ColumnWidth = 8.83
$ws.Range("Q13").Select()

$wb.Save()
